$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$n = $t.Rows.Count

# Simple single-run cell replacements (1-indexed rows)
$cell = $t.Cell(1, 1)
$cell.Range.Text = "0M"

$cell = $t.Cell(2, 1)
$cell.Range.Text = "0M"

$cell = $t.Cell(3, 1)
$cell.Range.Text = "0M"

$cell = $t.Cell(4, 1)
$cell.Range.Text = "202"

$cell = $t.Cell(5, 1)
$cell.Range.Text = "0.00002"

$cell = $t.Cell(6, 1)
$cell.Range.Text = "0.00008"

$cell = $t.Cell(9, 1)
$cell.Range.Text = "0.00004"

$cell = $t.Cell(10, 1)
$cell.Range.Text = "0.00004"

$cell = $t.Cell(11, 1)
$cell.Range.Text = "0.00004"

$cell = $t.Cell(12, 1)
$cell.Range.Text = "0.00707"

# Collapse the last three multi-run rows down to a single summary value each
$cell = $t.Cell($n - 2, 1)
$cell.Range.Text = "100"

$cell = $t.Cell($n - 1, 1)
$cell.Range.Text = "0.01"

$cell = $t.Cell($n, 1)
$cell.Range.Text = "166"
